# Add "lat" / "lon" columns (with latitude & longitude values) to the
# funicular station table, per:
#   "Aggiunte latitudine e longitudine fermate e stazioni"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Expand the Tabella1 table from A1:D5 to A1:F5 (adds two new, blank,
# auto-named columns which we then rename/populate below).
$lo.Resize($ws.Range("A1:F5")) | Out-Null

# Stage a template cell with the number format + alignment the new data
# columns should end up with, then paste just the formatting onto the
# target cells. Doing it this way (one shot via PasteSpecial) keeps the
# generated style table minimal, the same way Excel folds a single
# "format painter" operation into one cellXfs record instead of one per
# property assignment.
$tmpl = $ws.Range("H1")
$tmpl.NumberFormat = "0.00000000"
$tmpl.HorizontalAlignment = -4108  # xlCenter
$tmpl.VerticalAlignment = -4108    # xlCenter
$tmpl.WrapText = $true

# Header cells E1/F1: match the style already used by the other table
# headers (e.g. D1 "province").
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E1").Value2 = "lat"
$ws.Range("F1").Value2 = "lon"

# Data cells E2:F5: apply the staged numeric style, then clean up the
# scratch template cell.
$tmpl.Copy()
$ws.Range("E2:F5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$tmpl.Clear()

# Station coordinates.
$ws.Range("E2").Value2 = 45.7007227
$ws.Range("F2").Value2 = 9.6650854000000006
$ws.Range("E3").Value2 = 45.703040399999999
$ws.Range("F3").Value2 = 9.6651491000000007
$ws.Range("E4").Value2 = 45.7065403
$ws.Range("F4").Value2 = 9.6580069000000002
$ws.Range("E5").Value2 = 45.703340500000003
$ws.Range("F5").Value2 = 9.6645737999999994

# Best-fit-ish column widths for the new columns.
$ws.Columns.Item(5).ColumnWidth = 20
$ws.Columns.Item(6).ColumnWidth = 19

# Leave the selection where the author ended up.
$ws.Range("F7").Select() | Out-Null
